$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6051.5
$ws.Range("I86").Value = 6051.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6051.5
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -4928.5

$ws.Range("H88").Value = 1986752.2
$ws.Range("I88").Value = 4334.3335
$ws.Range("J88").Value = 3176203
$ws.Range("K88").Value = 4334.3335
$ws.Range("L88").Value = 3176203
$ws.Range("M88").Value = -3928.3335
$ws.Range("N88").Value = -3177015

$ws.Range("H89").Value = 6051.5
$ws.Range("I89").Value = 6051.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 30257.5
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -24641.5

$ws.Range("H91").Value = 1986752.2
$ws.Range("I91").Value = 4334.3335
$ws.Range("J91").Value = 3176203
$ws.Range("K91").Value = 4334.3335
$ws.Range("L91").Value = 3176203
$ws.Range("M91").Value = -2930.3335
$ws.Range("N91").Value = -3179011

$ws.Range("H129").Value = 996.1892
$ws.Range("I129").Value = 411.125
$ws.Range("J129").Value = 1157.5862
$ws.Range("K129").Value = 1233.375
$ws.Range("L129").Value = 3472.7586
$ws.Range("M129").Value = 3766.625
$ws.Range("N129").Value = -13472.7586

$ws.Range("H137").Value = 1145.5454
$ws.Range("I137").Value = 971.7143
$ws.Range("J137").Value = 1449.75
$ws.Range("K137").Value = 2915.1429
$ws.Range("L137").Value = 4349.25
$ws.Range("M137").Value = -365.1428999999998
$ws.Range("N137").Value = -9449.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8945.933999999999
$ws.Range("I32").Value = 3922.7407
$ws.Range("J32").Value = 54154.668
$ws.Range("K32").Value = 3922.7407
$ws.Range("L32").Value = 54154.668
$ws.Range("M32").Value = -3635.7407
$ws.Range("N32").Value = -54728.668

$ws.Range("H88").Value = 3215.1428
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3215.1428
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 3215.1428
$ws.Range("N88").Value = -4027.1428

$ws.Range("H91").Value = 3215.1428
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3215.1428
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 3215.1428
$ws.Range("N91").Value = -6023.1428

$ws.Range("H102").Value = 4341
$ws.Range("I102").Value = 3608.889
$ws.Range("J102").Value = 4999.9
$ws.Range("K102").Value = 3608.889
$ws.Range("L102").Value = 4999.9
$ws.Range("M102").Value = -1986.889
$ws.Range("N102").Value = -8243.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 50319.75
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 50319.75
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 50319.75
$ws.Range("N135").Value = -60459.75

$ws.Range("H137").Value = 49281.25
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 49281.25
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 49281.25
$ws.Range("N137").Value = -59481.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2169.6428
$ws.Range("I31").Value = 2053.5652
$ws.Range("J31").Value = 2703.6
$ws.Range("K31").Value = 2053.5652
$ws.Range("L31").Value = 2703.6
$ws.Range("M31").Value = -1758.5652
$ws.Range("N31").Value = -3293.6

$ws.Range("H34").Value = 2169.6428
$ws.Range("I34").Value = 2053.5652
$ws.Range("J34").Value = 2703.6
$ws.Range("K34").Value = 2053.5652
$ws.Range("L34").Value = 2703.6
$ws.Range("M34").Value = -1851.5652
$ws.Range("N34").Value = -3107.6

$ws.Range("H58").Value = 5764.095
$ws.Range("I58").Value = 844.5263
$ws.Range("J58").Value = 52500
$ws.Range("K58").Value = 844.5263
$ws.Range("L58").Value = 52500
$ws.Range("M58").Value = -641.5263
$ws.Range("N58").Value = -52906

$ws.Range("H136").Value = 5764.095
$ws.Range("I136").Value = 844.5263
$ws.Range("J136").Value = 52500
$ws.Range("K136").Value = 2533.5789
$ws.Range("L136").Value = 157500
$ws.Range("M136").Value = 16.42110000000002
$ws.Range("N136").Value = -162600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 51.82353
$ws.Range("I12").Value = 24.2
$ws.Range("J12").Value = 63.333332
$ws.Range("K12").Value = 72.59999999999999
$ws.Range("L12").Value = 189.999996
$ws.Range("M12").Value = 100.4
$ws.Range("N12").Value = -535.999996

$ws.Range("H68").Value = 1166.6666
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1166.6666
$ws.Range("K68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("M68").Value = 3499.9998
$ws.Range("N68").Value = -5121.9998

$ws.Range("H71").Value = 1166.6666
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1166.6666
$ws.Range("K71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("M71").Value = 10499.9994
$ws.Range("N71").Value = -18611.9994

$ws.Range("H131").Value = 14921.703
$ws.Range("I131").Value = 71894.21000000001
$ws.Range("J131").Value = 1628.1167
$ws.Range("K131").Value = 215682.63
$ws.Range("L131").Value = 4884.3501
$ws.Range("M131").Value = -210642.63
$ws.Range("N131").Value = -14964.3501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5123.6924
$ws.Range("I70").Value = 4418.909
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4418.909
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4148.909
$ws.Range("N70").Value = -9540

$ws.Range("H73").Value = 5123.6924
$ws.Range("I73").Value = 4418.909
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4418.909
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3482.909
$ws.Range("N73").Value = -10872

$ws.Range("H97").Value = 1854.8
$ws.Range("I97").Value = 2314.2856
$ws.Range("J97").Value = 1452.75
$ws.Range("K97").Value = 2314.2856
$ws.Range("L97").Value = 1452.75
$ws.Range("M97").Value = -1818.2856
$ws.Range("N97").Value = -2444.75

$ws.Range("H124").Value = 38495
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 38495
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 38495
$ws.Range("N124").Value = -48315

$ws.Range("H132").Value = 2774.3125
$ws.Range("I132").Value = 2534.15
$ws.Range("J132").Value = 3174.5833
$ws.Range("K132").Value = 7602.450000000001
$ws.Range("L132").Value = 9523.749899999999
$ws.Range("M132").Value = -5072.450000000001
$ws.Range("N132").Value = -14583.7499

$ws.Range("H136").Value = 22326
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 22326
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 66978
$ws.Range("N136").Value = -72078

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4762
$ws.Range("I122").Value = 9359.637000000001
$ws.Range("J122").Value = 2888.889
$ws.Range("K122").Value = 28078.911
$ws.Range("L122").Value = 8666.667000000001
$ws.Range("M122").Value = -25628.911
$ws.Range("N122").Value = -13566.667

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2068.4285
$ws.Range("I81").Value = 1913.1666
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 3826.3332
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -2765.3332
$ws.Range("N81").Value = -8122

$ws.Range("H84").Value = 2068.4285
$ws.Range("I84").Value = 1913.1666
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 19131.666
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -13827.666
$ws.Range("N84").Value = -40608

$ws.Range("H138").Value = 52330
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 52330
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 52330
$ws.Range("N138").Value = -62610
